# Update MSME definitions: replace literal "<br/>" markers with real line
# breaks inside the cell text values on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("B36").Value = "<50 Industry and Construction " + $nl + "<25 Agriculture" + $nl + "<15 Wholesale trade" + $nl + "<10 Retail trade, transportation, service, and other economic activities"

$ws.Range("D36").Value = "<500,00 thousand manat (454,5 thousand Euro) Industry and Construction " + $nl + "<250,00 thousand manat (227,3 thousand Euro) Agriculture" + $nl + "<1 Millionlion manat (900,9 thousand Euro) Wholesale trade" + $nl + "<250,00 thousand manat (227,3 thousand Euro) Retail trade, transportatio"

$ws.Range("B37").Value = ">50 Industry and Construction " + $nl + ">25 Agriculture" + $nl + ">15 Wholesale trade" + $nl + ">10 Retail trade, transportation, service, and other economic activities"

$ws.Range("D37").Value = ">=500,00 thousand manat Industry and Construction " + $nl + ">=250,00 thousand manat Agriculture" + $nl + ">=1 Millionlion manat Wholesale trade" + $nl + ">=250,00 thousand manat Retail trade, transportation, service, and other economic activities"
